$wb = $excel.ActiveWorkbook

# The "想去人数" (number of people interested) counts were refreshed for
# three rows on both the "展览" sheet and the "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 3308
    $ws.Range("F3").Value = 14
    $ws.Range("F5").Value = 1318
}
